# Region VI_QRF.xlsx edit
# - Delete the 5 "Negros Occidental" rows (216-220), shifting all rows
#   below them up by 5 (so old row 221 becomes new row 216, ... old row
#   228 becomes new row 223). This also naturally updates the sheet
#   <dimension> and the AA data-validation sqref.
# - Rename header cell B1 from "REGION" to "Region".
# - Narrow column C from width 19 to width 17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Delete rows 216 through 220 (Negros Occidental entries removed).
$ws.Range("A216:A220").EntireRow.Delete() | Out-Null

# Header text tweak.
$ws.Range("B1").Value = "Region"

# Column C width: stored xlsx width 19 -> 17.
# Excel's COM ColumnWidth property is offset from the raw OOXML column
# width by the default font padding (~0.8333 chars), so we compensate
# to land exactly on the target stored width.
$ws.Columns.Item(3).ColumnWidth = 17 - 0.8333333333333333
